# Apply cryptos list update (Tue May 21 16:41:41 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to remain plain text (matches original inlineStr cells),
    # so numeric-looking strings like "1.00" or "0.500" are not coerced to numbers
    # (which would drop the trailing zeros / introduce float artifacts).
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "69.899.68"
Set-TextValue $ws.Range("E2") "  +2.93%  "

Set-TextValue $ws.Range("D3") "3.797.66"
Set-TextValue $ws.Range("E3") "  +21.37%  "

Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  +0.36%  "

Set-TextValue $ws.Range("D5") "616.01"
Set-TextValue $ws.Range("E5") "  +6.64%  "

Set-TextValue $ws.Range("D6") "177.64"
Set-TextValue $ws.Range("E6") "  -1.31%  "

Set-TextValue $ws.Range("D7") "3.794.73"
Set-TextValue $ws.Range("E7") "  +21.60%  "

Set-TextValue $ws.Range("D8") "1.00"
Set-TextValue $ws.Range("E8") "  +0.18%  "

Set-TextValue $ws.Range("D9") "0.543"
Set-TextValue $ws.Range("E9") "  +4.82%  "

Set-TextValue $ws.Range("E10") "  +10.12%  "

Set-TextValue $ws.Range("D11") "6.37"
Set-TextValue $ws.Range("E11") "  -2.39%  "

Set-TextValue $ws.Range("D12") "0.500"
Set-TextValue $ws.Range("E12") "  +6.44%  "

Set-TextValue $ws.Range("D13") "40.59"
Set-TextValue $ws.Range("E13") "  +10.32%  "

Set-TextValue $ws.Range("D14") "0.0000258"
Set-TextValue $ws.Range("E14") "  +6.71%  "

Set-TextValue $ws.Range("D15") "4.455.90"
Set-TextValue $ws.Range("E15") "  +22.26%  "

Set-TextValue $ws.Range("D16") "3.807.82"
Set-TextValue $ws.Range("E16") "  +22.00%  "

Set-TextValue $ws.Range("D17") "70.201.11"
Set-TextValue $ws.Range("E17") "  +3.71%  "

Set-TextValue $ws.Range("D19") "7.57"
Set-TextValue $ws.Range("E19") "  +7.25%  "

Set-TextValue $ws.Range("D20") "518.17"
Set-TextValue $ws.Range("E20") "  +6.35%  "

Set-TextValue $ws.Range("D21") "16.68"
Set-TextValue $ws.Range("E21") "  +1.50%  "

Set-TextValue $ws.Range("D22") "9.49"
Set-TextValue $ws.Range("E22") "  +22.40%  "

Set-TextValue $ws.Range("D23") "0.739"
Set-TextValue $ws.Range("E23") "  +6.60%  "

Set-TextValue $ws.Range("D24") "88.96"
Set-TextValue $ws.Range("E24") "  +5.99%  "

Set-TextValue $ws.Range("D25") "2.48"
Set-TextValue $ws.Range("E25") "  +6.30%  "

Set-TextValue $ws.Range("D26") "13.58"
Set-TextValue $ws.Range("E26") "  +5.86%  "

Set-TextValue $ws.Range("D27") "10.88"
Set-TextValue $ws.Range("E27") "  +2.01%  "

Set-TextValue $ws.Range("D28") "0.0000127"
Set-TextValue $ws.Range("E28") "  +34.12%  "

Set-TextValue $ws.Range("E29") "  +0.04%  "

Set-TextValue $ws.Range("D30") "2.51"
Set-TextValue $ws.Range("E30") "  +6.76%  "

Set-TextValue $ws.Range("D31") "2.86"
Set-TextValue $ws.Range("E31") "  +9.08%  "

Set-TextValue $ws.Range("D32") "7.88"
Set-TextValue $ws.Range("E32") "  -3.20%  "

Set-TextValue $ws.Range("D33") "32.19"
Set-TextValue $ws.Range("E33") "  +13.99%  "

Set-TextValue $ws.Range("E34") "  +1.87%  "

Set-TextValue $ws.Range("E35") "  -0.11%  "

Set-TextValue $ws.Range("D36") "6.22"
Set-TextValue $ws.Range("E36") "  +10.70%  "

Set-TextValue $ws.Range("E37") "  +9.76%  "

Set-TextValue $ws.Range("D38") "0.343"
Set-TextValue $ws.Range("E38") "  +5.63%  "

Set-TextValue $ws.Range("D39") "2.18"
Set-TextValue $ws.Range("E39") "  +6.83%  "

Set-TextValue $ws.Range("E40") "  +6.48%  "

Set-TextValue $ws.Range("D41") "51.46"
Set-TextValue $ws.Range("E41") "  +4.48%  "

Set-TextValue $ws.Range("B42") "Cosmos"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D42") "8.86"
Set-TextValue $ws.Range("E42") "  +5.88%  "

Set-TextValue $ws.Range("B43") "Arweave"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextValue $ws.Range("D43") "44.50"
Set-TextValue $ws.Range("E43") "  -8.94%  "

Set-TextValue $ws.Range("D44") "425.28"
Set-TextValue $ws.Range("E44") "  +10.13%  "

Set-TextValue $ws.Range("D45") "3.118.95"
Set-TextValue $ws.Range("E45") "  +12.19%  "

Set-TextValue $ws.Range("D46") "2.70"
Set-TextValue $ws.Range("E46") "  +0.59%  "

Set-TextValue $ws.Range("E47") "  +4.99%  "

Set-TextValue $ws.Range("D48") "27.83"
Set-TextValue $ws.Range("E48") "  +3.66%  "

Set-TextValue $ws.Range("E49") "  +0.88%  "

Set-TextValue $ws.Range("D51") "2.47"
Set-TextValue $ws.Range("E51") "  +5.99%  "
